$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the static "NBL01xx" product-code values in B3:B8 with a shared
# formula that generates a random "PK_BD_#####" code. Setting the formula on
# the whole B3:B8 range at once makes the engine record it as a single
# shared formula group (t="shared"), matching row 3 as the "master" cell.
$ws.Range("B3:B8").Formula = '="PK_BD_"&TEXT(RANDBETWEEN(0,99999),"00000")'

# The old cells used a taller, wrapped row (two-line "NBLxxxx\n" text), the
# new single-line generated codes don't need that extra height, so let
# Excel re-fit rows 3:8 back down to the default height.
$ws.Rows("3:8").AutoFit()

# Add a new (currently empty) row below the table for the next entry, with
# the same left/top/wrap formatting used throughout the sheet's body font.
$c = $ws.Cells.Item(9, 2)
$c.Value = ""
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4160
$c.WrapText = $true
$c.Font.Name = "Arial"

# Move the active selection to where the user was last working.
$ws.Range("D14").Select()
